$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the surrounding double quotes from the schedule strings across
# the data rows (rows 2-4, columns B-G). Column H (CERRADO) is untouched.
$ws.Range("B2").Value = "9:00-18:00 "
$ws.Range("C2").Value = "12:00-21:00 "
$ws.Range("D2").Value = "9:00-18:00"
$ws.Range("E2").Value = "12:00-21:00"
$ws.Range("F2").Value = "9:00-18:00"
$ws.Range("G2").Value = "9:00-18:00"

$ws.Range("B3").Value = "9:00-18:00 "
$ws.Range("C3").Value = "12:00-21:00"
$ws.Range("D3").Value = "9:00-18:00"
$ws.Range("E3").Value = "12:00-21:00"
$ws.Range("F3").Value = "9:00-18:00"
$ws.Range("G3").Value = "9:00-12:00"

$ws.Range("B4").Value = "9:00-18:00 "
$ws.Range("C4").Value = "12:00-21:00 "
$ws.Range("D4").Value = "9:00-18:00"
$ws.Range("E4").Value = "12:00-21:00"
$ws.Range("F4").Value = "9:00-18:00"
$ws.Range("G4").Value = "9:00-12:00"

# Update the selected cell in the sheet view to H4
$ws.Range("H4").Select()
